# Fix Training Data Issue (#48)
# Column BF ("Date") held a malformed date string "6-27-2007-08" for every
# data row (the value was off by a day due to how the NBA stats source
# displayed dates). Re-stamp it as a plain ISO date string "2008-06-27".
#
# NumberFormat is forced to Text ("@") before the assignment so Excel does
# not silently reinterpret the ISO-looking string as a date serial, then
# restored to the default style afterwards so the cell's formatting is left
# untouched - only its value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "6-27-2007-08") {
        $cell.NumberFormat = "@"
        $cell.Value = "2008-06-27"
        $cell.Style = "Normal"
    }
}
